$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# Set values in the exact order the shared-strings table needs them added,
# so new unique strings end up appended in the same sequence as the target.
$ws.Range("A113").Value = "em_ui_global_cooldown"
$ws.Range("C113").Value = "Global Request Cooldown(s)"
$ws.Range("A114").Value = "em_ui_tab_whitelist"
$ws.Range("C114").Value = "ホワイトリスト"
$ws.Range("D114").Value = "白名单"
$ws.Range("D113").Value = "全局请求间隔(s)"
$ws.Range("A115").Value = "em_ui_whitelist"
$ws.Range("C115").Value = "AIサービスホワイトリストモード：{0}"
$ws.Range("D115").Value = "AI服务白名单模式: {0}"

# Update selection to D115
$ws.Range("D115").Select()

# Update workbook window position
$excel.ActiveWindow.Left = 10965
$excel.ActiveWindow.Top = 2460
